# PowerPointBlank.pptx was re-saved with the title-slide placeholders
# (Title 1 / Subtitle 2) removed from the one slide in the template,
# leaving an empty slide behind (the title-bearing variant moved to the
# new PowerPointWithTitle.pptx template instead).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove every shape on the slide (currently the ctrTitle "Title 1" and
# subTitle "Subtitle 2" placeholders), walking backwards so the live
# Shapes collection doesn't shift indices out from under us.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $s.Shapes.Item($i).Delete()
}
